# Generate Report for Handback
# Applies: rename overview.md -> overview-test2.md (status flips to "not in
# sync"), and appends a brand-new "authoringResource.md" row (status "in
# sync") to the Overview sheet and the zh-cn detail sheet.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276   # BGR for RGB(100,149,237) == FF6495ED

function Style-AsHyperlink($range) {
    $range.Font.Name = "Calibri"
    $range.Font.Underline = $true
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Existing row 2: file got renamed, so it is no longer in sync.
$wsOverview.Range("B2").Value = "Handed back: not in sync with en-US"

# New row 3 for authoringResource.md (in sync).
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"

# Rebuild the hyperlinks on this sheet (existing link needs new display
# text + new row needs a fresh link) in one pass so nothing is orphaned.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A2"),
    "https://github.com/OpenLocalizationOrg/PowerShell-Docs/blob/8c06757253ac5d54a77bb008b87bd426238fff1a/dsc/testmove/overview-test2.md",
    "", "", "overview-test2.md") | Out-Null
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A3"),
    "https://github.com/OpenLocalizationOrg/PowerShell-Docs/blob/8c06757253ac5d54a77bb008b87bd426238fff1a/dsc/testmove/authoringResource.md",
    "", "", "authoringResource.md") | Out-Null
Style-AsHyperlink $wsOverview.Range("A2")
Style-AsHyperlink $wsOverview.Range("A3")

# Extend the "Overview" table to cover the new row.
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:B3"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Existing row 2 updates (rename propagates into this row too).
$wsZhCn.Range("B2").Value = ".md"
$wsZhCn.Range("C2").Value = "Handed back: not in sync with en-US"
$wsZhCn.Range("E2").Value = "2016-04-12 06:50:20"
$wsZhCn.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H2").Value = "2016-04-12 16:11:00"
$wsZhCn.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("J2").Value = "Include"

# New row 3 for authoringResource.md.
$wsZhCn.Range("A3").Value = "authoringResource.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "authoringResource.106f1f4d77068b2b95f1927b611b5f405ec7a317.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-04-12 19:20:28"
$wsZhCn.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("F3").Value = "authoringResource.md"
$wsZhCn.Range("G3").Value = "authoringResource.106f1f4d77068b2b95f1927b611b5f405ec7a317.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-04-12 22:40:49"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("J3").Value = "Include"

# Rebuild all hyperlinks on this sheet together (existing A2/D2/F2/G2 plus
# the four new ones on row 3).
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A2"),
    "https://github.com/OpenLocalizationOrg/PowerShell-Docs/blob/8c06757253ac5d54a77bb008b87bd426238fff1a/dsc/testmove/overview-test2.md",
    "", "", "overview-test2.md") | Out-Null
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("D2"),
    "https://github.com/OpenLocalizationOrg/olhandoff/blob/293ea328e8c483d6c1a09d74eb64ff2bdd51198e/ol-handoff/OpenLocalizationOrg/PowerShell-Docs.zh-cn/master/high/overview.8cc85dd99239e10c76baa6006d906abfd6122c3f.zh-cn.xlf",
    "", "", "overview.8cc85dd99239e10c76baa6006d906abfd6122c3f.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationOrg/PowerShell-Docs.zh-cn/blob/967632d5fdfa88e703c6c5bb2a6b4c7f93fee227/dsc/testmove/overview-test2.md",
    "", "", "overview-test2.md") | Out-Null
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G2"),
    "https://github.com/OpenLocalizationOrg/olhandback/blob/d9570c2718d1e7bddb99bbaec358f36b800661ce/ol-handback/OpenLocalizationOrg/PowerShell-Docs.zh-cn/master/high/overview.8cc85dd99239e10c76baa6006d906abfd6122c3f.zh-cn.xlf",
    "", "", "overview.8cc85dd99239e10c76baa6006d906abfd6122c3f.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A3"),
    "https://github.com/OpenLocalizationOrg/PowerShell-Docs/blob/8c06757253ac5d54a77bb008b87bd426238fff1a/dsc/testmove/authoringResource.md",
    "", "", "authoringResource.md") | Out-Null
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("D3"),
    "https://github.com/OpenLocalizationOrg/olhandoff/blob/293ea328e8c483d6c1a09d74eb64ff2bdd51198e/ol-handoff/OpenLocalizationOrg/PowerShell-Docs.zh-cn/master/high/authoringResource.106f1f4d77068b2b95f1927b611b5f405ec7a317.zh-cn.xlf",
    "", "", "authoringResource.106f1f4d77068b2b95f1927b611b5f405ec7a317.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F3"),
    "https://github.com/OpenLocalizationOrg/PowerShell-Docs.zh-cn/blob/967632d5fdfa88e703c6c5bb2a6b4c7f93fee227/dsc/testmove/authoringResource.md",
    "", "", "authoringResource.md") | Out-Null
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G3"),
    "https://github.com/OpenLocalizationOrg/olhandback/blob/d9570c2718d1e7bddb99bbaec358f36b800661ce/ol-handback/OpenLocalizationOrg/PowerShell-Docs.zh-cn/master/high/authoringResource.106f1f4d77068b2b95f1927b611b5f405ec7a317.zh-cn.xlf",
    "", "", "authoringResource.106f1f4d77068b2b95f1927b611b5f405ec7a317.zh-cn.xlf") | Out-Null

foreach ($addr in @("A2", "D2", "F2", "G2", "A3", "D3", "F3", "G3")) {
    Style-AsHyperlink $wsZhCn.Range($addr)
}

# Re-apply the date/time number format to the correspond-datetime columns,
# since adding the hyperlinks above only touched column A/D/F/G.
$wsZhCn.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Extend the "zh-cn" table to cover the new row.
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:L3"))

Write-Output "done"
